# RDCC-5182 Added Version check
# Insert a new "VERSION" worksheet at the front of the workbook containing a
# small "File version" / "vx.xx" key-value pair in row 6, pushing the
# existing "Staff Data" and "Sheet2" tabs one position to the right.

$wb = $excel.ActiveWorkbook

$versionSheet = $wb.Worksheets.Add()
$versionSheet.Name = "VERSION"

$versionSheet.Range("A6").Value = "File version"
$versionSheet.Range("B6").Value = "vx.xx"

$versionSheet.Range("B6").Select()
